$d = $word.ActiveDocument

# ===========================================================================
# Helper: locate the paragraph index whose text contains a given substring.
# ===========================================================================
function Find-ParagraphIndex($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# ===========================================================================
# Change 1: insert a brand-new paragraph right after the paragraph that ends
# with "...così da rendere il programma più concorrente possibile." and
# right before the "Organizzazione" (Titolo3) heading.
#
# New paragraph text:
#   "La suddivisione in diversi " + italic("package, ") +
#   "invece, è dovuta ad una maggior divisione dei vari aspetti del
#    progetto, così da aumentare l'indipendenza delle varie funzionalità
#    del progetto."
# ===========================================================================
$anchorIdx = Find-ParagraphIndex("concorrente possibile")
$anchorPara = $d.Paragraphs.Item($anchorIdx)
$endPos = $anchorPara.Range.End
$insPoint = $d.Range($endPos - 1, $endPos - 1)
$insPoint.InsertParagraphAfter()

$newIdx = $anchorIdx + 1
$newPara = $d.Paragraphs.Item($newIdx)
$s = $newPara.Range.Start
$d.Range($s, $s).InsertAfter("La suddivisione in diversi ")

$newPara = $d.Paragraphs.Item($newIdx)
$e = $newPara.Range.End
$d.Range($e - 1, $e - 1).InsertAfter("package, ")
$italicStart = $e - 1
$italicEnd = $italicStart + [int]("package, ".Length)
$d.Range($italicStart, $italicEnd).Font.Italic = 1

$newPara = $d.Paragraphs.Item($newIdx)
$e = $newPara.Range.End
$d.Range($e - 1, $e - 1).InsertAfter("invece, è dovuta ad una maggior divisione dei vari aspetti del progetto, così da aumentare l'indipendenza delle varie funzionalità del progetto.")

# ===========================================================================
# Change 2: append bold "PuzzleParser." (with spell-check proof markers) to
# the end of the paragraph "Questo package contiene l'inerfaccia IParser e
# la sua implementazione ".
#
# We build this via InsertXML (which always creates a brand-new following
# paragraph) and then delete the paragraph mark that separates the two
# paragraphs, merging "PuzzleParser." back into the previous paragraph while
# preserving the w:proofErr spell-check markers that came from the XML.
# ===========================================================================
$implIdx = Find-ParagraphIndex("e la sua implementazione")
$implPara = $d.Paragraphs.Item($implIdx)
$e = $implPara.Range.End
$insRange = $d.Range($e - 1, $e - 1)

$puzzleParserXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>PuzzleParser</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insRange.InsertXML($puzzleParserXml)

# delete the paragraph mark that now separates "...implementazione " from
# "PuzzleParser." so the two merge back into a single paragraph
$implPara = $d.Paragraphs.Item($implIdx)
$markEnd = $implPara.Range.End
$d.Range($markEnd - 1, $markEnd).Delete()

# ===========================================================================
# Change 3: in the following paragraph (which used to start with
# "Il parser con la sua interfaccia sono posti nel package parser, così come
# la classe SequentSort ...") remove:
#   a) the leading "Il parser con la sua interfaccia sono posti nel package
#      parser, " (including the bold "parser" run and the bold ", " run)
#   b) the sentence "Ciò è dovuto ad una maggior divisione dei vari aspetti
#      del progetto, così da aumentare l'indipendenza delle varie
#      funzionalità del progetto. " in the middle of the paragraph.
# ===========================================================================
$d.Content.Find.Execute("Il parser con la sua interfaccia sono posti nel package parser, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

$d.Content.Find.Execute("Ciò è dovuto ad una maggior divisione dei vari aspetti del progetto, così da aumentare l'indipendenza delle varie funzionalità del progetto. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

Write-Host "All changes applied"
